# Change the first "Potential dispersal rate: additive" table-cell caption
# to "Potential dispersal rate: interactive".  The source run `": additive"`
# is split into `": "` (keeps the original run's rsidR) and a new run
# holding just the word "interactive", matching how Word represents a
# retyped word inside existing text.

$d = $word.ActiveDocument

# Locate the first (and, for this edit, only relevant) occurrence of the
# full caption text so we don't disturb the later rows/tables that already
# read "...interactive".
$rng = $d.Content
$found = $rng.Find.Execute("Potential dispersal rate: additive", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found -and $rng.Find.Found) {
    $paraRange = $d.Range($rng.Start, $rng.End)

    $newParagraphXml = '<w:p w14:paraId="2A4389B0" w14:textId="3AC26303" w:rsidR="001709BE" w:rsidRPr="000E27A9" w:rsidRDefault="001709BE" w:rsidP="002D1460"><w:pPr><w:keepNext/><w:spacing w:after="0" w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr></w:pPr><w:r w:rsidRPr="000E27A9"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Potential dispersal rate</w:t></w:r><w:r w:rsidR="001E646E"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>interactive</w:t></w:r></w:p>'

    $xmlWrapped = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>REPLACE_PARA</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

    $xmlWrapped = $xmlWrapped.Replace("REPLACE_PARA", $newParagraphXml)

    $paraRange.InsertXML($xmlWrapped)
}
